$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D6").Value = "2016-03-08 16:37:12"
$wsZh.Range("C7").Value = "e6786125-8ef7-4cb8-b17f-c35cb3dc22a1.7ef090da5caa9a544f29d239573985732162f22b.zh-cn.xlf"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D6").Value = "2016-03-08 16:37:18"
$wsDe.Range("C7").Value = "e6786125-8ef7-4cb8-b17f-c35cb3dc22a1.7ef090da5caa9a544f29d239573985732162f22b.de-de.xlf"
